$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, Coin(B), Link(C), Price(D), Volume1h(E), ForceTextOnD
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '68.266.01', '  +0.96%  ', $false),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.353.23', '  +0.49%  ', $false),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.999', '  +0.12%  ', $true),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '584.24', '  +0.46%  ', $true),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '177.10', '  +0.41%  ', $true),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.999', '  -0.12%  ', $true),
    @(8, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.591', '  +0.01%  ', $true),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.183', '  +2.08%  ', $true),
    @(10, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.581', '  +0.55%  ', $true),
    @(11, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '48.10', '  +5.71%  ', $true),
    @(12, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000274', '  +1.41%  ', $true),
    @(13, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '685.39', '  +3.04%  ', $true),
    @(14, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.894.40', '  +0.61%  ', $false),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '8.42', '  -0.17%  ', $true),
    @(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '68.300.12', '  +0.95%  ', $false),
    @(17, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.120', '  +1.17%  ', $true),
    @(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.351.29', '  +0.91%  ', $false),
    @(19, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '17.45', '  +0.12%  ', $true),
    @(20, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '11.21', '  +2.05%  ', $true),
    @(21, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.894', '  +0.35%  ', $true),
    @(22, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '5.46', '  -0.67%  ', $true),
    @(23, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '16.99', '  -0.66%  ', $true),
    @(24, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '100.45', '  +0.85%  ', $true),
    @(25, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.91', '  +1.34%  ', $true),
    @(26, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.70', '  +0.88%  ', $true),
    @(27, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '9.51', '  +2.15%  ', $true),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '33.00', '  -2.15%  ', $true),
    @(29, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '8.52', '  +0.75%  ', $true),
    @(30, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '6.94', '  -6.78%  ', $true),
    @(31, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.08', '  +0.83%  ', $true),
    @(32, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '555.64', '  -3.92%  ', $true),
    @(33, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.105', '  +0.84%  ', $true),
    @(34, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '58.03', '  +2.49%  ', $true),
    @(35, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  -0.03%  ', $true),
    @(36, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '3.716.12', '  +0.16%  ', $false),
    @(37, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '3.32', '  -2.04%  ', $true),
    @(38, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.138', '  +5.09%  ', $true),
    @(39, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '34.89', '  +1.03%  ', $true),
    @(40, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '3.17', '  +1.59%  ', $true),
    @(41, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '2.61', '  -0.60%  ', $true),
    @(42, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0673', '  +0.05%  ', $false),
    @(43, 'TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.336', '  +0.21%  ', $true),
    @(44, 'ApeXProtocol', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', '3.23', '  -1.73%  ', $true),
    @(45, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0411', '  +1.09%  ', $true),
    @(46, 'ThetaToken', 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta', '2.64', '  +1.82%  ', $true),
    @(47, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.128', '  +0.29%  ', $true),
    @(48, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  -0.05%  ', $true),
    @(49, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '1.34', '  -0.83%  ', $true),
    @(50, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '131.92', '  +2.92%  ', $true),
    @(51, 'CoreDAO', 'https://coinranking.com/coin/HFvoXUQh4+coredao-core', '2.57', '  -1.19%  ', $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[5]) {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Updated $($data.Count) rows"